# Revert "chore: remove view" — re-add the trial-version notice paragraph
# and the "fwqfq" paragraph that used to precede the section break.

$d = $word.ActiveDocument

# Start at the very beginning of the document body (before the sectPr).
$rng = $d.Content
$rng.Collapse(1)

# First paragraph: the red "trial version" notice.
$rng.InsertAfter("This document was generated by a trial version of Telerik Document Processing.")
$rng.InsertParagraphAfter()

# Second paragraph: plain text "fwqfq".
$rng.Collapse(0)
$rng.InsertAfter("fwqfq")

# Color only the first paragraph's run text red (FF0000), leaving the
# paragraph mark / second paragraph unaffected.
$p1 = $d.Paragraphs(1).Range
$p1Text = $d.Range($p1.Start, $p1.End - 1)
$p1Text.Font.Color = 255
